$wb = $excel.ActiveWorkbook

# zh-cn sheet: Correspond Handoff Datetime (D) and Correspond Handback DateTime (G)
# for rows 21 and 22 (5803e10a... and 62e2ac14... entries)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D21").Value = "2016-03-08 07:07:16"
$wsZh.Range("D22").Value = "2016-03-08 07:07:16"
$wsZh.Range("G21").Value = "2016-03-08 07:08:00"
$wsZh.Range("G22").Value = "2016-03-08 07:08:00"

# de-de sheet: same rows / columns
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D21").Value = "2016-03-08 07:07:27"
$wsDe.Range("D22").Value = "2016-03-08 07:07:27"
$wsDe.Range("G21").Value = "2016-03-08 07:08:17"
$wsDe.Range("G22").Value = "2016-03-08 07:08:17"
